$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy formatting (bold, centered, border) from E1 then set text
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamp values for F2:F13
$timestamps = @(
    "2021-10-05 10:52:04.530112",
    "2021-10-05 10:52:04.530123",
    "2021-10-05 10:52:04.530126",
    "2021-10-05 10:52:04.530129",
    "2021-10-05 10:52:04.530133",
    "2021-10-05 10:52:04.530136",
    "2021-10-05 10:52:04.530138",
    "2021-10-05 10:52:04.530141",
    "2021-10-05 10:52:04.530144",
    "2021-10-05 10:52:04.530147",
    "2021-10-05 10:52:04.530149",
    "2021-10-05 10:52:04.530152"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
